$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data for columns A (phone number, stored as text) and B (amount)
$data = @(
  @("9375204365",0),
  @("9375204382",0),
  @("9375204397",0),
  @("9375204518",0),
  @("9375204529",0),
  @("9375204526",0),
  @("9375204522",0),
  @("9375204537",0),
  @("9375204553",0),
  @("9375205286",0),
  @("9375205296",0),
  @("9375204585",0),
  @("9375205306",0),
  @("9375204605",0),
  @("9375204608",0),
  @("9375204602",0),
  @("9375205341",0),
  @("9375205312",0),
  @("9992653911",0),
  @("9376255745",0),
  @("9397467135",0),
  @("9372870453",0),
  @("9375823270",0),
  @("9372921150",0),
  @("9372920857",0),
  @("9377710154",0),
  @("9377710159",0),
  @("9372975879",0),
  @("9372921364",0),
  @("9375203244",0),
  @("9375814225",0),
  @("9397498073",0),
  @("9372919584",0),
  @("9376268702",0),
  @("9376258078",0),
  @("9375204625",0),
  @("9375204613",0),
  @("9375204623",0),
  @("9375204622",0),
  @("9375204619",0),
  @("9372914234",0),
  @("9375204654",0),
  @("9375204651",0),
  @("9375204649",0),
  @("9375204513",0),
  @("9372919472",0),
  @("9375204634",0),
  @("9375204632",0),
  @("9375204653",0),
  @("9372912146",0),
  @("9372912164",0),
  @("9372912104",0),
  @("9274559669",800),
  @("9274222297",800),
  @("9397431818",800),
  @("9297614242",800),
  @("9272403636",800),
  @("9274593737",800),
  @("9376285757",800),
  @("9274844222",800),
  @("9375933773",800),
  @("9274424554",800),
  @("9270429393",800),
  @("9376284343",800),
  @("9397465151",800),
  @("9297629191",800),
  @("9397432626",800),
  @("9397492121",800),
  @("9297635454",800),
  @("9397482525",800),
  @("9397455252",800),
  @("9397421515",800),
  @("9297629494",800),
  @("9297631414",800),
  @("9297624848",800),
  @("9297605252",800),
  @("9297624949",800),
  @("9297615959",800),
  @("9397468282",800),
  @("9376293434",800),
  @("9372946767",800),
  @("9274986161",800),
  @("9397459393",800),
  @("9397466262",800),
  @("9397452626",800),
  @("9397468181",800),
  @("9397462525",800),
  @("9397455353",800),
  @("9397458282",800),
  @("9376295454",800),
  @("9376295858",800),
  @("9376298181",800),
  @("9376293131",800),
  @("9376294747",800),
  @("9397409696",800),
  @("9274716969",800),
  @("9377715522",4000),
  @("9397454884",4000),
  @("9274062662",4000),
  @("9375852244",4000),
  @("9274758855",4000),
  @("9375751991",4000),
  @("9376188877",4000),
  @("9274827227",4000),
  @("9274826886",4000),
  @("9272495454",4000),
  @("9376199977",4000),
  @("9274606161",4000),
  @("9274559559",4000),
  @("9274363737",4000),
  @("9397407700",4000),
  @("9397447700",4000),
  @("9376211177",4000),
  @("9274888843",4000),
  @("9370011115",4000),
  @("9375937799",4000),
  @("9376168558",4000),
  @("9375979559",4000),
  @("9272436262",4000),
  @("9270373377",4000),
  @("9370079966",4000),
  @("9370076699",4000),
  @("9376202266",4000),
  @("9375881155",4000),
  @("9375909955",4000),
  @("9274686644",4000),
  @("9375272255",4000),
  @("9274266644",4000),
  @("9376199966",4000),
  @("9397472277",4000),
  @("9274948844",4000),
  @("9274229944",4000),
  @("9274845544",4000),
  @("9376299966",4000),
  @("9274549944",4000),
  @("9397313377",4000),
  @("9370051100",4000),
  @("9370065500",4000),
  @("9272466622",4000),
  @("9274373344",4000),
  @("9274858833",4000),
  @("9372976666",5000),
  @("9375726666",5000),
  @("9375826666",5000),
  @("9372936666",5000),
  @("9375736666",15000),
  @("9274796666",5000),
  @("9375816666",5000),
  @("9274286666",5000),
  @("9274596666",5000),
  @("9274396666",5000),
  @("9375296666",5000),
  @("9274386666",5000),
  @("9272444446",15000),
  @("9274888885",15000)
)

$rowCount = $data.Length

# Force column A to text formatting so phone numbers are written as strings
# (shared-string cells), matching the source list, not as numbers.
$ws.Range("A1:A$rowCount").NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Drop the temporary text formatting again so cells fall back to the default
# (unstyled) look, while keeping the values stored as text.
$ws.Cells.ClearFormats()
